$d = $word.ActiveDocument

# Locate the paragraph that currently holds the lone "_GoBack" bookmark - that is
# the "...calc(100vh- (hh+hf))." paragraph this edit appends new content after.
$bm = $d.Bookmarks.Item("_GoBack")
$bmParaIndex = $bm.Range.Paragraphs.First.Index

$startPara = $d.Paragraphs.Item($bmParaIndex)
$endPara = $d.Paragraphs.Item($d.Paragraphs.Count)

# Range spanning from the start of that paragraph through to the end of the
# document (the trailing empty paragraph) - InsertXML below replaces this whole
# span with the paragraphs reconstructed from it plus the new content.
$rng = $d.Range($startPara.Range.Start, $endPara.Range.End)

# Paragraph 1 (existing, edited): close the "...calc(100vh- (hh+hf))." paragraph
# right after its final run - the bookmark that used to sit here moves down into
# the new "Зроблено ..." paragraph below.
$p_calc = '<w:p w:rsidR="00EA0A68" w:rsidRPr="0078160F" w:rsidRDefault="0078160F" w:rsidP="00BB1086"><w:pPr><w:rPr><w:lang w:val="ru-RU"/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">Вирішено: переробити хедер без </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>аватару</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, і на блок з проектами встановити висоту через </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>calc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="0078160F"><w:rPr><w:lang w:val="ru-RU"/></w:rPr><w:t>(100</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>vh</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="0078160F"><w:rPr><w:lang w:val="ru-RU"/></w:rPr><w:t>- (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>hh</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="ru-RU"/></w:rPr><w:t>+</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>hf</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="0078160F"><w:rPr><w:lang w:val="ru-RU"/></w:rPr><w:t>))</w:t></w:r><w:r><w:rPr><w:lang w:val="ru-RU"/></w:rPr><w:t>.</w:t></w:r></w:p>'

# Paragraph 2 (new): blank spacer line.
$p_blank1 = '<w:p><w:pPr><w:rPr><w:lang w:val="ru-RU"/></w:rPr></w:pPr></w:p>'

# Paragraph 3 (new): "Додано на хедер та футер:" lead-in line.
$p_dodano = '<w:p><w:pPr><w:rPr><w:lang w:val="ru-RU"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="ru-RU"/></w:rPr><w:t>Додано на хедер та футер:</w:t></w:r></w:p>'

# Paragraph 4 (new): the syntax-highlighted "box-shadow: 0 0 10px 5px rgba(104, 102, 102, 0.8);" code line.
$p_codeblock = '<w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="212121"/><w:spacing w:after="0" w:line="345" w:lineRule="atLeast"/><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="EEFFFF"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:val="en-US" w:eastAsia="uk-UA"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="B2CCD6"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:eastAsia="uk-UA"/></w:rPr><w:t>box-shadow</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="89DDFF"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:eastAsia="uk-UA"/></w:rPr><w:t>:</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="EEFFFF"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:eastAsia="uk-UA"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="F78C6C"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:eastAsia="uk-UA"/></w:rPr><w:t>0</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="EEFFFF"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:eastAsia="uk-UA"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="F78C6C"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:eastAsia="uk-UA"/></w:rPr><w:t>0</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="EEFFFF"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:eastAsia="uk-UA"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="F78C6C"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:eastAsia="uk-UA"/></w:rPr><w:t>10px</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="EEFFFF"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:eastAsia="uk-UA"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="F78C6C"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:eastAsia="uk-UA"/></w:rPr><w:t>5px</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="EEFFFF"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:eastAsia="uk-UA"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="82AAFF"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:eastAsia="uk-UA"/></w:rPr><w:t>rgba</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="89DDFF"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:eastAsia="uk-UA"/></w:rPr><w:t>(</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="F78C6C"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:eastAsia="uk-UA"/></w:rPr><w:t>104</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="89DDFF"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:eastAsia="uk-UA"/></w:rPr><w:t>,</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="EEFFFF"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:eastAsia="uk-UA"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="F78C6C"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:eastAsia="uk-UA"/></w:rPr><w:t>102</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="89DDFF"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:eastAsia="uk-UA"/></w:rPr><w:t>,</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="EEFFFF"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:eastAsia="uk-UA"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="F78C6C"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:eastAsia="uk-UA"/></w:rPr><w:t>102</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="89DDFF"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:eastAsia="uk-UA"/></w:rPr><w:t>,</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="EEFFFF"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:eastAsia="uk-UA"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="F78C6C"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:eastAsia="uk-UA"/></w:rPr><w:t>0.8</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="89DDFF"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:eastAsia="uk-UA"/></w:rPr><w:t>)</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="89DDFF"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:val="en-US" w:eastAsia="uk-UA"/></w:rPr><w:t>;</w:t></w:r></w:p>'

# Paragraph 5 (new): blank spacer line.
$p_blank2 = '<w:p/>'

# Paragraph 6 (new): "Зроблено адаптивність ..." summary line - this is where the
# "_GoBack" bookmark now lives, in between its two runs.
$p_zrobleno = '<w:p><w:r><w:t>Зроблено</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve"> адаптивність сторінки проектів, додано функціонал відображення та приховання навігації</w:t></w:r></w:p>'

# Paragraph 7 (existing, unchanged): the trailing empty paragraph that was
# already the last paragraph of the body before the edit.
$p_trailing = '<w:p w:rsidR="00CB5A24" w:rsidRPr="00667D1A" w:rsidRDefault="00CB5A24" w:rsidP="00BB1086"/>'

$frag = $p_calc + $p_blank1 + $p_dodano + $p_codeblock + $p_blank2 + $p_zrobleno + $p_trailing

$rng.InsertXML($frag)
